# #5: property boat&car done
# Sheet "汽車" (car) was missing its header row and the extra metadata
# columns that every other sheet in this workbook carries: row 1 had
# literal data instead of column headers, and row 2 stopped at column G.
# This fills in the proper header row (name, capacity, owner,
# register_date, register_reason, acquire_value, property_category,
# category, date, legislator_name, legislator_id, source_file, index) and
# completes the single data row to match the rest of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Header row (row 1): proper column names instead of duplicated data ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2): fill in the remaining columns ---
$ws.Range("A2").Value = 31
$ws.Range("B2").Value = "LEXUS"
$ws.Range("C2").Value = 4608
$ws.Range("D2").Value = "田正超"
$ws.Range("E2").Value = "97年07月07日"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = 2100000
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# Force this one to stay plain text - otherwise an ISO-looking
# "yyyy-mm-dd" string gets auto-parsed into a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-03-26"
$ws.Range("K2").Value = "潘維剛"
$ws.Range("L2").Value = 678
$ws.Range("M2").Value = "tmp71a01"
$ws.Range("N2").Value = 31

# Match formatting: header cells look like B1 (bold / bordered / centered),
# data cells look like B2 - including J2, so the text-number-format tweak
# above doesn't leave it looking different from its row neighbours.
$ws.Range("B1").Copy()
$ws.Range("C1:N1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
